# Updated cryptos list on Tue Mar  7 23:26:46 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set the Price (D) cell as plain text, even when the text looks
# like a number (Excel would otherwise silently convert it to a float and
# change its stored type). Prefixing with an apostrophe forces text entry;
# resetting Style afterwards removes the "quote prefix" cell style that the
# apostrophe trick would otherwise leave behind, so the cell's style stays
# the same as before (no explicit style / default).
function Set-Price($row, $value) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

function Set-Volume($row, $value) {
    $ws.Range("E$row").Value = $value
}

function Set-Row($row, $d, $e) {
    Set-Price $row $d
    Set-Volume $row $e
}

Set-Row 2  "22.129.00"     "  -1.47%  "
Set-Row 3  "1.555.17"      "  -1.06%  "
Set-Row 4  "0.9989"        "  -0.24%  "
Set-Row 5  "0.9997"        "  -0.15%  "
Set-Row 6  "288.12"        "  -0.14%  "
Set-Row 7  "0.3795"        "  +2.36%  "
Set-Volume 8 "  -0.70%  "
Set-Row 9  "43.20"         "  -10.65%  "
Set-Row 10 "1.139"         "  +0.39%  "
Set-Row 11 "0.07361"       "  -1.83%  "
Set-Row 12 "0.9991"        "  -0.28%  "
Set-Row 13 "20.17"         "  -2.72%  "
Set-Row 14 "5.823"         "  -1.85%  "
Set-Row 15 "6.822"         "  -0.72%  "
Set-Row 16 "1.558.97"      "  -0.61%  "
Set-Row 17 "0.00001103"    "  -1.48%  "
Set-Row 18 "0.06616"       "  -1.91%  "
Set-Row 19 "85.88"         "  -2.00%  "
Set-Row 20 "6.393"         "  +0.64%  "
Set-Row 21 "0.9994"        "  -0.15%  "
Set-Row 22 "16.11"         "  -2.58%  "
Set-Row 23 "11.70"         "  -2.89%  "
Set-Row 24 "22.101.62"     "  -1.60%  "
Set-Row 25 "2.312"         "  -3.35%  "
Set-Row 26 "2.532"         "  -1.82%  "
Set-Row 27 "150.64"        "  -1.94%  "
Set-Volume 28 "  -2.94%  "
Set-Row 29 "4.911"         "  -2.17%  "
Set-Row 30 "121.65"        "  -2.30%  "
Set-Row 31 "1.732.17"      "  -0.79%  "
Set-Row 32 "1.082"         "  +2.14%  "
Set-Row 33 "5.971"         "  -2.41%  "
Set-Row 34 "1.858"         "  -7.66%  "
Set-Row 35 "0.08222"       "  -1.78%  "
Set-Row 36 "9.340"         "  -4.64%  "
Set-Row 37 "0.02337"       "  -5.41%  "

# Rows 38 and 39 swap places: Hedera now ranks above
# InternetComputer(DFINITY), while the rank numbers in column A are
# unchanged.
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Row 38 "0.06256" "  -2.50%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-Row 39 "5.295" "  -0.91%  "

Set-Volume 40 "  -4.59%  "
Set-Row 41 "1.255"         "  -2.70%  "
Set-Row 42 "11.05"         "  -2.25%  "
Set-Row 43 "0.6059"        "  -4.10%  "
Set-Row 44 "0.9993"        "  -0.17%  "
Set-Row 45 "13.76"         "  -0.28%  "
Set-Row 46 "3.737"         "  -0.93%  "
Set-Row 47 "0.5858"        "  -5.04%  "
Set-Row 48 "1.990"         "  -3.47%  "
Set-Row 49 "122.37"        "  -2.83%  "
Set-Row 50 "1.177"         "  -3.02%  "
Set-Row 51 "0.07018"       "  -2.83%  "
